$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws.Range("D2").Value = 3830.39
$ws.Range("E2").Value = -3830.39

$ws.Range("D4").Value = 4378.85
$ws.Range("E4").Value = 13121.15
$ws.Range("F4").Value = 0.25022
